# Update the Tgfb2/Tgfbr1 LR-pair table per Dr Hou's advice.
# Sending/target cluster now iterate over ECs/FAPs/sCs (x) ECs/FAPs/M2/sCs,
# and all the expression-derived statistics are recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb2"
$ws.Range("C2").Value = "Tgfbr1"
$ws.Range("D2").Value = "ECs"
$values = @(2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 44.50020533333333, 133.500616, 0.2926972930209797, 0.2926972930209797, 58.97773896906044, 530.799650721544, 0.00854113299644775, 0.00854113299644775)
$cols = @("E2", "F2", "G2", "H2", "I2", "J2", "K2", "L2", "M2", "N2", "O2", "P2", "Q2", "R2", "S2", "T2")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb2"
$ws.Range("C3").Value = "Tgfbr1"
$ws.Range("D3").Value = "FAPs"
$values = @(2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 28.185334, 84.55600199999999, 0.1853872561462678, 0.1853872561462678, 37.35504721733533, 336.195424956018, 0.005409743268374896, 0.005409743268374895)
$cols = @("E3", "F3", "G3", "H3", "I3", "J3", "K3", "L3", "M3", "N3", "O3", "P3", "Q3", "R3", "S3", "T3")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfb2"
$ws.Range("C4").Value = "Tgfbr1"
$ws.Range("D4").Value = "M2"
$values = @(2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 62.31760433333333, 186.952813, 0.4098901108273345, 0.4098901108273344, 82.59178522925743, 743.3260670633169, 0.01196090990241593, 0.01196090990241592)
$cols = @("E4", "F4", "G4", "H4", "I4", "J4", "K4", "L4", "M4", "N4", "O4", "P4", "Q4", "R4", "S4", "T4")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tgfb2"
$ws.Range("C5").Value = "Tgfbr1"
$ws.Range("D5").Value = "sCs"
$values = @(2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 17.031762, 51.09528599999999, 0.1120253400054181, 0.1120253400054181, 22.57281299928599, 203.155316993574, 0.003268985914024058, 0.003268985914024058)
$cols = @("E5", "F5", "G5", "H5", "I5", "J5", "K5", "L5", "M5", "N5", "O5", "P5", "Q5", "R5", "S5", "T5")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tgfb2"
$ws.Range("C6").Value = "Tgfbr1"
$ws.Range("D6").Value = "ECs"
$values = @(3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 44.50020533333333, 133.500616, 0.2926972930209797, 0.2926972930209797, 910.2451288963457, 8192.206160067111, 0.1318213421058903, 0.1318213421058903)
$cols = @("E6", "F6", "G6", "H6", "I6", "J6", "K6", "L6", "M6", "N6", "O6", "P6", "Q6", "R6", "S6", "T6")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tgfb2"
$ws.Range("C7").Value = "Tgfbr1"
$ws.Range("D7").Value = "FAPs"
$values = @(3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 28.185334, 84.55600199999999, 0.1853872561462678, 0.1853872561462678, 576.5268449356793, 5188.741604421113, 0.08349239127666905, 0.08349239127666903)
$cols = @("E7", "F7", "G7", "H7", "I7", "J7", "K7", "L7", "M7", "N7", "O7", "P7", "Q7", "R7", "S7", "T7")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Tgfb2"
$ws.Range("C8").Value = "Tgfbr1"
$ws.Range("D8").Value = "M2"
$values = @(3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 62.31760433333333, 186.952813, 0.4098901108273345, 0.4098901108273344, 1274.697394405427, 11472.27654964884, 0.1846011760734612, 0.1846011760734612)
$cols = @("E8", "F8", "G8", "H8", "I8", "J8", "K8", "L8", "M8", "N8", "O8", "P8", "Q8", "R8", "S8", "T8")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Tgfb2"
$ws.Range("C9").Value = "Tgfbr1"
$ws.Range("D9").Value = "sCs"
$values = @(3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 17.031762, 51.09528599999999, 0.1120253400054181, 0.1120253400054181, 348.3821766864779, 3135.439590178301, 0.05045257001513991, 0.0504525700151399)
$cols = @("E9", "F9", "G9", "H9", "I9", "J9", "K9", "L9", "M9", "N9", "O9", "P9", "Q9", "R9", "S9", "T9")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tgfb2"
$ws.Range("C10").Value = "Tgfbr1"
$ws.Range("D10").Value = "ECs"
$values = @(3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 44.50020533333333, 133.500616, 0.2926972930209797, 0.2926972930209797, 1051.893598992263, 9467.042390930368, 0.1523348179186416, 0.1523348179186416)
$cols = @("E10", "F10", "G10", "H10", "I10", "J10", "K10", "L10", "M10", "N10", "O10", "P10", "Q10", "R10", "S10", "T10")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Tgfb2"
$ws.Range("C11").Value = "Tgfbr1"
$ws.Range("D11").Value = "FAPs"
$values = @(3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 28.185334, 84.55600199999999, 0.1853872561462678, 0.1853872561462678, 666.2434970350773, 5996.191473315695, 0.0964851216012239, 0.09648512160122388)
$cols = @("E11", "F11", "G11", "H11", "I11", "J11", "K11", "L11", "M11", "N11", "O11", "P11", "Q11", "R11", "S11", "T11")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Tgfb2"
$ws.Range("C12").Value = "Tgfbr1"
$ws.Range("D12").Value = "M2"
$values = @(3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 62.31760433333333, 186.952813, 0.4098901108273345, 0.4098901108273344, 1473.060373806047, 13257.54336425442, 0.2133280248514573, 0.2133280248514573)
$cols = @("E12", "F12", "G12", "H12", "I12", "J12", "K12", "L12", "M12", "N12", "O12", "P12", "Q12", "R12", "S12", "T12")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Tgfb2"
$ws.Range("C13").Value = "Tgfbr1"
$ws.Range("D13").Value = "sCs"
$values = @(3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 17.031762, 51.09528599999999, 0.1120253400054181, 0.1120253400054181, 402.595927213392, 3623.363344920528, 0.05830378407625413, 0.05830378407625413)
$cols = @("E13", "F13", "G13", "H13", "I13", "J13", "K13", "L13", "M13", "N13", "O13", "P13", "Q13", "R13", "S13", "T13")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $values[$i]
}
